$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.416.01"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "3.492.18"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.02"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.65"
$ws.Range("E6").Value = "  +4.17%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +1.75%  "

$ws.Range("D9").Value = "3.494.34"
$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  +4.88%  "

$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").Value = "4.096.94"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.34"
$ws.Range("E14").Value = "  +10.33%  "

$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").Value = "67.418.60"
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "3.490.38"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.31"
$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.37"
$ws.Range("E21").Value = "  -1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.94"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.11"
$ws.Range("E23").Value = "  +1.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.542"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000122"
$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.35"
$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("E29").Value = "  -2.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.22"
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("E33").Value = "  +0.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.54"
$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("E35").Value = "  +0.69%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.83"
$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  +10.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.89"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.79"
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.64"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").Value = "2.847.08"
$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.36"
$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.06"
$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.71"
$ws.Range("E48").Value = "  -2.28%  "

$ws.Range("E49").Value = "  -0.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "334.32"
$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("E51").Value = "  -1.36%  "
